# change from dictionary to list and on progress of choose the best path
# Update the "address" column (D) values on the "Child" sheet to the new
# coordinate pairs.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Child")

$values = @{
    2  = "4,1"
    3  = "-9,0"
    4  = "0,-5"
    5  = "0,-2"
    6  = "9,-1"
    7  = "-5,0"
    8  = "-10,-9"
    9  = "9,3"
    10 = "9,3"
    11 = "-9,-3"
    12 = "-5,8"
    13 = "-2,9"
    14 = "-3,9"
    15 = "-1,1"
    16 = "6,0"
    17 = "1,8"
    18 = "5,4"
    19 = "0,4"
    20 = "2,6"
    21 = "2,-3"
    22 = "3,8"
}

foreach ($row in $values.Keys) {
    $ws.Range("D$row").Value = $values[$row]
}
